# Admin can create single staff - simplify the staff roster sheet by removing
# the "nickname" column (previously column D) and highlighting the remaining
# data with a solid white fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "nickname" column (column D). All the columns to the
# right (Email, Phone) shift left to fill the gap, and the sheet's used
# range/dimension shrinks from A1:F10 to A1:E10 automatically.
$ws.Columns.Item(4).Delete()

# Apply a solid white background fill across the remaining data range.
$ws.Range("A1:E10").Interior.Color = 16777215
